$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "add select all option" - the last document (index 4, row 6) is dropped
# from the TF-IDF matrix and the remaining documents' term weights are
# recomputed, so the whole (now shorter) table can be selected/used as-is.
$ws.Rows("6:6").Delete()

# Recomputed TF-IDF weights for the remaining documents (row 2 = doc 0 ... row 5 = doc 3)
$ws.Range("B2").Value = 0.9853284724047701
$ws.Range("C2").Value = 0.1342617650160634
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0.02048060822278933
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0.004551246271730962
$ws.Range("I2").Value = 0.002275623135865481
$ws.Range("J2").Value = 0.03868559330971318
$ws.Range("K2").Value = 0.06144182466836799
$ws.Range("L2").Value = 0.006826869407596443
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0.02226731011084226
$ws.Range("P2").Value = 0.06958534409638206

$ws.Range("B3").Value = 0.1077132271348133
$ws.Range("C3").Value = 0.1320939316800721
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0.6604696584003606
$ws.Range("G3").Value = 0.1995708967262472
$ws.Range("H3").Value = 0.1761252422400962
$ws.Range("I3").Value = 0.04403131056002405
$ws.Range("J3").Value = 0.08806262112004809
$ws.Range("K3").Value = 0.6604696584003606
$ws.Range("L3").Value = 0.04403131056002405
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0.06652363224208238
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0.1077132271348133

$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 0.108722157306214
$ws.Range("D4").Value = 0.3285204422586898
$ws.Range("E4").Value = 0.2083432623952618
$ws.Range("F4").Value = 0.108722157306214
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0.05436107865310699
$ws.Range("I4").Value = 0.05436107865310699
$ws.Range("J4").Value = 0.108722157306214
$ws.Range("K4").Value = 0.8697772584497118
$ws.Range("L4").Value = 0.05436107865310699
$ws.Range("M4").Value = 0.2083432623952618
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0.06649140279703235
$ws.Range("P4").Value = 0

$ws.Range("B5").Value = 0.1110330316926266
$ws.Range("C5").Value = 0.1072816477638703
$ws.Range("D5").Value = 0.01246798941549944
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0.3631071155084841
$ws.Range("G5").Value = 0.1371478835704938
$ws.Range("H5").Value = 0.00825243444337464
$ws.Range("I5").Value = 0.01650486888674928
$ws.Range("J5").Value = 0.2145632955277406
$ws.Range("K5").Value = 0.6106801488097233
$ws.Range("L5").Value = 0.03300973777349856
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 0.01246798941549944
$ws.Range("O5").Value = 0.6359164542395885
$ws.Range("P5").Value = 0.02018782394411392

